$d = $word.ActiveDocument

# --- Edit 1: "(BRYAN)" -> "(BRYAN" + " TO SEND" + ")" (three separate runs) ---
$f1 = $d.Content
$f1.Find.Execute("(BRYAN)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f1.Find.Found) {
    $r1 = $d.Range($f1.Start, $f1.End)
    $xml1 = '<w:p><w:r><w:t>(BRYAN</w:t></w:r><w:r><w:t xml:space="preserve"> TO SEND</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>'
    $r1.InsertXML($xml1)
}

# --- Edit 2: append a new run "(who will?)" right after "WILL ABSORB THIS COST" ---
$f2 = $d.Content
$f2.Find.Execute("WILL ABSORB THIS COST", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f2.Find.Found) {
    $r2 = $d.Range($f2.End, $f2.End)
    $r2.InsertAfter("(who will?)")
}
